$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1506605198818927
$ws.Range("C2").Value = 0.9984471741463956
$ws.Range("D2").Value = 0.3276022338459857
$ws.Range("F2").Value = "Pipeline(steps=[('model', AdaBoostRegressor(n_estimators=150))])"
$ws.Range("G2").Value = 0.1217136106832186
$ws.Range("H2").Value = 0.9740000000000001
